$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHART of ACCOUNTS_Updated")

$ws.Range("H1").Value = "Medical"
$ws.Range("I1").Value = "Dental"
$ws.Range("J1").Value = "Vision"
$ws.Range("K1").Value = "Life"
$ws.Range("G1").Value = "MEDICAL WAIVER"

$ws.Range("T4").Select()
